# "update number of bearings required"
# The F695 2RS Bearing line item's quantity note changes from x28 to x32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "x32"

# Leave the active cell on B4, matching the saved selection state.
$ws.Range("B4").Select()
